# Sheldon et al extraction for acclimation analysis
# Inserts a new literature-review row (row 11) for "Sheldon2020" ahead of
# the existing "Peng201432" row, shifting the rows below it down by one,
# and updates the dependent ranges (dimension, AutoFilter, _FilterDatabase,
# selection) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 11, pushing existing rows 11-16 to 12-17 ---
$ws.Rows(11).Insert()
$ws.Rows(11).RowHeight = 18

# --- Populate the new row with the Sheldon et al. 2020 record ---
$ws.Range("A11").Value = "Sheldon2020"
$ws.Range("B11").Value = "Sheldon, K.S. and Padash, M. and Carter, A.W. and Marshall, K.E."
$ws.Range("C11").Value = "Different amplitudes of temperature fluctuation induce distinct transcriptomic and metabolomic responses in the dung beetle Phanaeus vindex"
$ws.Range("D11").Value = "Journal of Experimental Biology"
$ws.Range("E11").Value = "10.1242/jeb.233239"
$ws.Range("F11").Value = 2020
$ws.Range("G11").Value = 'Most studies exploring molecular and physiological responses to temperature have focused on constant temperature treatments. To gain a better understanding of the impact of fluctuating temperatures, we investigated the effects of increased temperature variation on Phanaeus vindex dung beetles across levels of biological organization. Specifically, we hypothesized that increased temperature variation is energetically demanding. We predicted that thermal sensitivity of metabolic rate and energetic reserves would be reduced with increasing fluctuation. To test this, we examined the responses of dung beetles to constant (20C), low fluctuation (205C), or high fluctuation (2012C) temperature treatments using respirometry, assessment of energetic reserves and HPLC-MS-based metabolomics. We found no significant differences in metabolic rate or energetic reserves, suggesting increased fluctuations were not energetically demanding. To understand why there was no effect of increased amplitude of temperature fluctuation on energetics, we assembled and annotated a de novo transcriptome, finding non-overlapping transcriptomic and metabolomic responses of beetles exposed to different fluctuations. We found that 58 metabolites increased in abundance in both fluctuation treatments, but 15 only did so in response to high-amplitude fluctuations. We found that 120 transcripts were significantly upregulated following acclimation to any fluctuation, but 174 were upregulated only in beetles from the high-amplitude fluctuation treatment. Several differentially expressed transcripts were associated with post-translational modifications to histones that support a more open chromatin structure. Our results demonstrate that acclimation to different temperature fluctuations is distinct and may be supported by increasing transcriptional plasticity. Our results indicate for the first time that histone modifications may underlie rapid acclimation to temperature variation.  2020 Company of Biologists Ltd. All rights reserved.'
$ws.Range("H11").Value = "selected"
$ws.Range("I11").Value = "y"
$ws.Range("J11").Value = "y"
$ws.Range("K11").Value = "fig 1"
$ws.Range("M11").Value = "full analysis for figure 2"
$ws.Range("N11").Value = "y"
$ws.Range("O11").Value = "SCOPUS405"

# --- Re-apply the AutoFilter over the now-larger data range ---
$ws.AutoFilterMode = $false | Out-Null
$ws.Range("A2:T20").AutoFilter() | Out-Null

# --- Keep the cached _FilterDatabase defined name in sync ---
foreach ($n in $wb.Names) {
    if ($n.Name() -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$2:`$T`$20"
    }
}

# --- Match the recorded selection/active cell ---
$ws.Range("J22").Select() | Out-Null
